$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09838033333333333
$ws.Range("H2").Value = 0.295141
$ws.Range("I2").Value = 0.0257774858695505
$ws.Range("J2").Value = 0.0257774858695505
$ws.Range("M2").Value = 0.4049933333333333
$ws.Range("N2").Value = 1.21498
$ws.Range("Q2").Value = 0.03984337913111111
$ws.Range("R2").Value = 0.35859041218
$ws.Range("S2").Value = 0.0257774858695505
$ws.Range("T2").Value = 0.0257774858695505

# Row 3 (FAPs -> ECs)
$ws.Range("I3").Value = 0.4417479616037814
$ws.Range("J3").Value = 0.4417479616037814
$ws.Range("M3").Value = 0.4049933333333333
$ws.Range("N3").Value = 1.21498
$ws.Range("Q3").Value = 0.6827947303955555
$ws.Range("R3").Value = 6.14515257356
$ws.Range("S3").Value = 0.4417479616037814
$ws.Range("T3").Value = 0.4417479616037814

# Row 4 (Inflammatory-Mac -> ECs)
$ws.Range("G4").Value = 0.243138
$ws.Range("H4").Value = 0.729414
$ws.Range("I4").Value = 0.06370669977418356
$ws.Range("J4").Value = 0.06370669977418356
$ws.Range("M4").Value = 0.4049933333333333
$ws.Range("N4").Value = 1.21498
$ws.Range("Q4").Value = 0.09846926907999999
$ws.Range("R4").Value = 0.88622342172
$ws.Range("S4").Value = 0.06370669977418356
$ws.Range("T4").Value = 0.06370669977418356

# Row 5 (MuSCs -> ECs)
$ws.Range("G5").Value = 0.6601003333333334
$ws.Range("H5").Value = 1.980301
$ws.Range("I5").Value = 0.1729586233188772
$ws.Range("J5").Value = 0.1729586233188772
$ws.Range("M5").Value = 0.4049933333333333
$ws.Range("N5").Value = 1.21498
$ws.Range("Q5").Value = 0.2673362343311111
$ws.Range("R5").Value = 2.40602610898
$ws.Range("S5").Value = 0.1729586233188772
$ws.Range("T5").Value = 0.1729586233188772

# Row 6 (Resolving-Mac -> ECs)
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.128962333333333
$ws.Range("H6").Value = 3.386887
$ws.Range("I6").Value = 0.2958092294336073
$ws.Range("J6").Value = 0.2958092294336073
$ws.Range("M6").Value = 0.4049933333333333
$ws.Range("N6").Value = 1.21498
$ws.Range("Q6").Value = 0.4572222185844445
$ws.Range("R6").Value = 4.11499996726
$ws.Range("S6").Value = 0.2958092294336073
$ws.Range("T6").Value = 0.2958092294336073

$wb.Save()
